$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtered save games) for rows 2-6, columns B-G

$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 10.35301142835362

$ws.Range("B3").Value = 0.06328177979961902
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 71517.89157740913
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 71522.69069507453

$ws.Range("B4").Value = 1.505614041169197
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 0.7127328510149897
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 4.371470058157054

$ws.Range("B5").Value = 0.06328177979961902
$ws.Range("C5").Value = 0.05231270169004087
$ws.Range("D5").Value = 16.98373111632243
$ws.Range("E5").Value = 6.48142807727062
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 23.58075367508271

$ws.Range("B6").Value = 1.505614041169197
$ws.Range("C6").Value = 1.65323645889881
$ws.Range("D6").Value = 0.1529057820181812
$ws.Range("E6").Value = 0.4998867070740569
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 3.811642989160245
